# Auto-update draw results: append the 2025-11-22 Pick 4 draw as a new
# row (67) at the bottom of the "Results" sheet, mirroring the layout of
# every existing row (Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67

# Leading "'" forces Excel to keep these values as plain text instead of
# auto-converting them to a date serial / number, matching how every
# other row in the column is stored (t="str").
$ws.Cells.Item($newRow, 1).Value = "'2025-11-22"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "'251122"
$ws.Cells.Item($newRow, 4).Value = "7-6-9-3"
$ws.Cells.Item($newRow, 5).Value = "2025-11-22T21:36:33.961+04:00"
